$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty F3 cell entirely
$ws.Range("F3").ClearContents()

# Correct District names (column G) to official names
$ws.Range("G4").Value = "Chamarajanagar"
$ws.Range("G6").Value = "Kodagu"
$ws.Range("G14").Value = "Kodagu"
$ws.Range("G16").Value = "Chamarajanagar"
$ws.Range("G17").Value = "Chamarajanagar"
$ws.Range("G27").Value = "Kodagu"
$ws.Range("G28").Value = "Chamarajanagar"
$ws.Range("G31").Value = "Chamarajanagar"
